$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the data rows (2-6) with the column-name placeholder values,
# leaving the header row (row 1) untouched.
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = "sku"
    $ws.Cells.Item($row, 2).Value = "name"
    $ws.Cells.Item($row, 3).Value = "quantity"
    $ws.Cells.Item($row, 4).Value = "cost_per"
    $ws.Cells.Item($row, 5).Value = "total_cost"
}
